$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Administrator'
$ws.Range("G3").Value = 'Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Asmaa Reda, Administrator'
$ws.Range("G4").Value = 'Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda'
$ws.Range("G5").Value = 'Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Eman Tantawi'
$ws.Range("G6").Value = 'Dr. Alshimaa Atef, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat'
$ws.Range("G7").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Menna tu''Alllah Mohammad, Dr. Nada Mohammad, Dr. Fatma Elhady'
$ws.Range("G11").Value = 'Dr. Amal Awwad, Dr. Safa Hany, Dr. Aya Saeed'
$ws.Range("G12").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Marina Youhanna'
$ws.Range("G17").Value = 'Dr. Esraa Samy, Dr. Mohammad Safwat'
$ws.Range("G19").Value = 'Dr. Rania Ahmad Youssef, Dr. Mariam Toma Gerges'
$ws.Range("G20").Value = 'Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range("G25").Value = 'Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil'
$ws.Range("G28").Value = 'Dr. Maryam Ashraf, Dr. Aya Emad'
$ws.Range("G30").Value = 'Dr. Yassmen Ahmad, Dr. Shorok Mohammad, Dr. Aya Hanafy, Dr. Wafaa Ebida'
